$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 - Team Norris/Drumm
$ws.Range("B2").Value = 5.5
$ws.Range("C2").Value = 3.5
$ws.Range("E2").Value = 4.5
$ws.Range("F2").Value = 5
$ws.Range("G2").Value = 3

# Row 3 - Team Brian Drumm
$ws.Range("B3").Value = 4
$ws.Range("F3").Value = 4.5
$ws.Range("G3").Value = 3.5

# Row 4 - Team Liz
$ws.Range("B4").Value = 5.5
$ws.Range("C4").Value = 4.5

# Row 5 - Team Marty Drumm
$ws.Range("B5").Value = 4.5
$ws.Range("E5").Value = 6
$ws.Range("F5").Value = 3.5
$ws.Range("G5").Value = 2

# Update the active selection to H3 as recorded in the saved view state
$ws.Range("H3").Select()
